$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Lower Right Cell" value for the Variables List Block (D15)
# from E238 to E239 (block extended to include S_SLUIS).
$ws.Range("D15").Value = "E239"

# Update the active selection to reflect where the edit was made.
$ws.Range("D15").Select()
